# Adds two new worksheets ("AssetDataModel" and "Material") describing the
# new binary asset/material file formats, documenting the asset loader and
# the generic asset input/output stream mentioned in the commit message.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: AssetDataModel
# ---------------------------------------------------------------------
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$assetSheet = $wb.Worksheets.Add($null, $last)
$assetSheet.Name = "AssetDataModel"

$assetSheet.Range("B3").Value = "<8 byte> "
$assetSheet.Range("C3").Value = "Magic Number"
$assetSheet.Range("D3").Value = '"VALASSET"'

$assetSheet.Range("B4").Value = "<1 byte>"
$assetSheet.Range("C4").Value = "Number of entries"

$assetSheet.Range("C6").Value = "<8 byte> TYPE_ID"
$assetSheet.Range("C7").Value = "<32 byte> Name"
$assetSheet.Range("C8").Value = "<4 byte> Offset"

# Best-fit the two populated columns, same as Excel's "AutoFit Column Width"
$assetSheet.Columns.Item(2).AutoFit() | Out-Null
$assetSheet.Columns.Item(3).AutoFit() | Out-Null

$assetSheet.Range("C9").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet 2: Material
# ---------------------------------------------------------------------
$last2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$materialSheet = $wb.Worksheets.Add($null, $last2)
$materialSheet.Name = "Material"

$materialSheet.Range("B3").Value = "Shader graph"

$materialSheet.Range("B5").Value = "CLASS:"
$materialSheet.Range("C5").Value = "vkShaderGraphAssetLoader"

$materialSheet.Range("B8").Value = "<uint16>"
$materialSheet.Range("C8").Value = "Num Nodes"

$materialSheet.Range("C9").Value = "<uint32>"
$materialSheet.Range("D9").Value = "Local IDX"

$materialSheet.Range("C10").Value = "<string>"
$materialSheet.Range("D10").Value = "Node Class Name"

$materialSheet.Range("C11").Value = "<vec2>"
$materialSheet.Range("D11").Value = "Location within Editor"

$materialSheet.Range("B13").Value = "<uint16> "
$materialSheet.Range("C13").Value = "Num Inputs"

$materialSheet.Range("C14").Value = "<uint32>"
$materialSheet.Range("D14").Value = "Local Node Input IDX"

$materialSheet.Range("C15").Value = "<uint8>"
$materialSheet.Range("D15").Value = "Input IDX"

$materialSheet.Range("C16").Value = "<uint8>"
$materialSheet.Range("D16").Value = "Input type (0 = Const, 1=ExternalNode)"

$materialSheet.Range("C17").Value = "<float>"
$materialSheet.Range("D17").Value = "Const float"

$materialSheet.Range("C18").Value = "<uint32>"
$materialSheet.Range("D18").Value = "Local OutputNode IDX"

$materialSheet.Range("C19").Value = "<uint8> "
$materialSheet.Range("D19").Value = "Output IDX"

$materialSheet.Range("B21").Value = "<uint16> "
$materialSheet.Range("C21").Value = "Num Graph Inputs"

$materialSheet.Range("C22").Value = "<string>"
$materialSheet.Range("D22").Value = "Attribute"

$materialSheet.Range("C23").Value = "<uint32>"
$materialSheet.Range("D23").Value = "Local OutputNode IDX"

$materialSheet.Range("C24").Value = "<uint8> "
$materialSheet.Range("D24").Value = "Output IDX"

$materialSheet.Range("B26").Value = "<uint16> "
$materialSheet.Range("C26").Value = "Num Graph Attributes"

$materialSheet.Range("C27").Value = "<string> "
$materialSheet.Range("D27").Value = "Attribute"

$materialSheet.Range("C28").Value = "<types>"

$materialSheet.Range("C28").Select() | Out-Null
